$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.406.22"
$ws.Range("E2").Value = "  +0.80%  "

$ws.Range("D3").Value = "3.624.12"
$ws.Range("E3").Value = "  +2.77%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.89%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "196.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.01%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.626"
$ws.Range("D7").Style = "Normal"

$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.213"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.53%  "

$ws.Range("E10").Value = "  -0.58%  "

$ws.Range("E11").Value = "  -1.10%  "

$ws.Range("E12").Value = "  +0.76%  "

$ws.Range("E13").Value = "  +0.30%  "

$ws.Range("D14").Value = "4.199.24"
$ws.Range("E14").Value = "  +2.82%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "606.03"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.26%  "

$ws.Range("E16").Value = "  +1.24%  "

$ws.Range("D17").Value = "70.476.79"
$ws.Range("E17").Value = "  +0.64%  "

$ws.Range("D18").Value = "3.623.13"
$ws.Range("E18").Value = "  +2.81%  "

$ws.Range("E19").Value = "  -0.35%  "

$ws.Range("E20").Value = "  +1.32%  "

$ws.Range("E21").Value = "  +0.42%  "

$ws.Range("E22").Value = "  -0.80%  "

$ws.Range("E23").Value = "  -1.24%  "

$ws.Range("E24").Value = "  +1.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.55%  "

$ws.Range("E26").Value = "  -6.42%  "

$ws.Range("E27").Value = "  -2.70%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.71%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.81"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.23%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.64"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.93%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.99%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.33%  "

$ws.Range("E33").Value = "  +0.68%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.24"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.19%  "

$ws.Range("D35").Value = "0.0₃0885"
$ws.Range("E35").Value = "  +3.81%  "

$ws.Range("D36").Value = "3.945.93"
$ws.Range("E36").Value = "  +5.57%  "

$ws.Range("E37").Value = "  +0.17%  "

$ws.Range("E38").Value = "  -0.51%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "516.66"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.93%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.69"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.13%  "

$ws.Range("E41").Value = "  -1.26%  "

$ws.Range("E42").Value = "  -2.38%  "

$ws.Range("E43").Value = "  +2.34%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0461"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.45%  "

$ws.Range("E45").Value = "  +6.66%  "

$ws.Range("E46").Value = "  +2.03%  "

$ws.Range("E47").Value = "  -0.18%  "

$ws.Range("E48").Value = "  +0.31%  "

$ws.Range("E49").Value = "  -0.28%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000249"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.20%  "

$ws.Range("E51").Value = "  +0.58%  "
